# Update the answers in the two-digit-division worksheet table.
# The table has 5 "data" rows (1, 5, 9, 13, 17) each holding 5 answer
# cells; the rows in between are blank spacer rows used for students'
# work. Several old values repeat verbatim elsewhere in the table
# (e.g. "48÷5=9, 3" appears twice), so a blind Find/Replace across the
# whole document could hit the wrong occurrence. Addressing each cell
# directly by (row, column) guarantees the correct cell is updated
# while leaving its run formatting (font/size) untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "29÷7=4, 1"
$t.Cell(1, 2).Range.Text  = "52÷8=6, 4"
$t.Cell(1, 3).Range.Text  = "74÷4=18, 2"
$t.Cell(1, 4).Range.Text  = "55÷3=18, 1"
$t.Cell(1, 5).Range.Text  = "64÷4=16, 0"

$t.Cell(5, 1).Range.Text  = "84÷4=21, 0"
$t.Cell(5, 2).Range.Text  = "29÷2=14, 1"
$t.Cell(5, 3).Range.Text  = "83÷2=41, 1"
$t.Cell(5, 4).Range.Text  = "75÷7=10, 5"
$t.Cell(5, 5).Range.Text  = "50÷3=16, 2"

$t.Cell(9, 1).Range.Text  = "22÷8=2, 6"
$t.Cell(9, 2).Range.Text  = "27÷9=3, 0"
$t.Cell(9, 3).Range.Text  = "57÷3=19, 0"
$t.Cell(9, 4).Range.Text  = "64÷2=32, 0"
$t.Cell(9, 5).Range.Text  = "88÷4=22, 0"

$t.Cell(13, 1).Range.Text = "57÷5=11, 2"
$t.Cell(13, 2).Range.Text = "62÷8=7, 6"
$t.Cell(13, 3).Range.Text = "71÷9=7, 8"
$t.Cell(13, 4).Range.Text = "34÷7=4, 6"
$t.Cell(13, 5).Range.Text = "31÷7=4, 3"

$t.Cell(17, 1).Range.Text = "52÷4=13, 0"
$t.Cell(17, 2).Range.Text = "27÷8=3, 3"
$t.Cell(17, 3).Range.Text = "48÷5=9, 3"
$t.Cell(17, 4).Range.Text = "47÷9=5, 2"
$t.Cell(17, 5).Range.Text = "73÷5=14, 3"
